$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.292.62'
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").Value = '3.468.90'
$ws.Range("E3").Value = '  +2.51%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.29%  '

$ws.Range("D7").Value = '3.467.95'
$ws.Range("E7").Value = '  +2.52%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.477'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.69'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.04%  '

$ws.Range("E11").Value = '  +0.92%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.392'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.14%  '

$ws.Range("D13").Value = '4.060.37'
$ws.Range("E13").Value = '  +2.47%  '

$ws.Range("E14").Value = '  -0.93%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.14%  '

$ws.Range("E16").Value = '  +0.71%  '

$ws.Range("D17").Value = '3.462.95'
$ws.Range("E17").Value = '  +2.22%  '

$ws.Range("D18").Value = '62.260.27'
$ws.Range("E18").Value = '  +0.62%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.24'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.566'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.62%  '

$ws.Range("D24").Value = '3.591.85'
$ws.Range("E24").Value = '  +1.98%  '

$ws.Range("B25").Value = 'LEO'
$ws.Range("C25").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '72.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.83%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("E28").Value = '  +1.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.180'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.87'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.65%  '

$ws.Range("E32").Value = '  +0.55%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.19'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.19%  '

$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '24.09'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.99%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.07'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.42%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.58'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.62%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '166.81'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0798'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.58%  '

$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.796'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.70%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.97'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.47%  '

$ws.Range("E44").Value = '  +1.17%  '

$ws.Range("E45").Value = '  -0.17%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.26'
$ws.Range("D46").Style = "Normal"

$ws.Range("E47").Value = '  +2.21%  '

$ws.Range("E48").Value = '  +0.37%  '

$ws.Range("D49").Value = '2.653.19'
$ws.Range("E49").Value = '  +11.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.87'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.30%  '

$ws.Range("E51").Value = '  +0.65%  '
